# Applies updated "想去人数" (want-to-go count) figures to the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 69
$wsExhibition.Range("F4").Value = 480
$wsExhibition.Range("F5").Value = 4657
$wsExhibition.Range("F6").Value = 363
$wsExhibition.Range("F8").Value = 286
$wsExhibition.Range("F9").Value = 720
$wsExhibition.Range("F10").Value = 202

$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F3").Value = 1

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F3").Value = 69
$wsAllTypes.Range("F4").Value = 480
$wsAllTypes.Range("F5").Value = 4657
$wsAllTypes.Range("F6").Value = 363
$wsAllTypes.Range("F8").Value = 286
$wsAllTypes.Range("F9").Value = 720
$wsAllTypes.Range("F11").Value = 202
$wsAllTypes.Range("F12").Value = 1
